$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 for the Chemours Chambers Works / Deepwater facility,
# shifting the old rows 9 (Islechem LLC) and 10 (Chemours - Corpus Christi Plant) down to 10 and 11.
$ws.Rows.Item(9).Insert()

# --- Update facility Label (col A) and City (col B) text, now in Title Case, for the rows that were not inserted ---
$ws.Cells.Item(2,1).Value = "Daikin America Inc."
$ws.Cells.Item(2,2).Value = "Decatur"
$ws.Cells.Item(3,1).Value = "Chemours El Dorado"
$ws.Cells.Item(3,2).Value = "El Dorado"
$ws.Cells.Item(4,1).Value = "Chemours Louisville Works"
$ws.Cells.Item(4,2).Value = "Louisville"
$ws.Cells.Item(5,1).Value = "Iofina Chemical Inc."
$ws.Cells.Item(5,2).Value = "Covington"
$ws.Cells.Item(6,1).Value = "ARKEMA, INC."
$ws.Cells.Item(6,2).Value = "Calvert City"
$ws.Cells.Item(7,1).Value = "Honeywell International - Geismar Complex"
$ws.Cells.Item(7,2).Value = "Geismar"
$ws.Cells.Item(8,1).Value = "Mexichem Fluor Inc."
$ws.Cells.Item(8,2).Value = "Saint Gabriel"
$ws.Cells.Item(10,1).Value = "Islechem LLC"
$ws.Cells.Item(10,2).Value = "Grand Island"
$ws.Cells.Item(11,1).Value = "Chemours - Corpus Christi Plant"
$ws.Cells.Item(11,2).Value = "Gregory"

# --- Convert GHG_co2e (col E) and pop_sq_mile_1mi (col G) from text to numeric values for the same rows ---
$ws.Cells.Item(2,5).Value = 20.5752959055157
$ws.Cells.Item(2,7).Value = 142.549590220656
$ws.Cells.Item(3,5).Value = 103.612502034919
$ws.Cells.Item(3,7).Value = 45.6605130373706
$ws.Cells.Item(4,5).Value = 36.0627684203524
$ws.Cells.Item(4,7).Value = 166.099283620707
$ws.Cells.Item(5,5).Value = 11.777227783969
$ws.Cells.Item(5,7).Value = 1405.25430123102
$ws.Cells.Item(6,5).Value = 74.6569158825578
$ws.Cells.Item(6,7).Value = 53.5114523922272
$ws.Cells.Item(7,5).Value = 120.261973778333
$ws.Cells.Item(7,7).Value = 75.7429777162124
$ws.Cells.Item(8,5).Value = 20.8371989696875
$ws.Cells.Item(8,7).Value = 79.5212447896903
$ws.Cells.Item(10,5).Value = 17.2846381134759
$ws.Cells.Item(10,7).Value = 419.968295103647
$ws.Cells.Item(11,5).Value = 130.239383628461
$ws.Cells.Item(11,7).Value = 40.0646859239259

# --- Fix blank/numeric column C (GHG blockgroups raw) for shifted rows ---
# Row 10 (Islechem LLC) should now be blank in column C
$ws.Cells.Item(10,3).ClearContents()
# Row 11 (Chemours - Corpus Christi Plant) should carry the original 17240 value
$ws.Cells.Item(11,3).Value = 17240

# --- Populate the new row 9: Chemours Chambers Works / Deepwater ---
$ws.Cells.Item(9,1).Value = "Chemours Chambers Works"
$ws.Cells.Item(9,2).Value = "Deepwater"
$ws.Cells.Item(9,3).Value = 2619
$ws.Cells.Item(9,4).Value = 5
$ws.Cells.Item(9,5).Value = 19.373996565504
$ws.Cells.Item(9,6).Value = 5935
$ws.Cells.Item(9,7).Value = 306.338445964601
$ws.Cells.Item(9,8).Value = 0
$ws.Cells.Item(9,9).Value = 0.8
$ws.Cells.Item(9,10).Value = 5156
$ws.Cells.Item(9,11).Value = 390
$ws.Cells.Item(9,12).Value = 30
$ws.Cells.Item(9,13).Value = 70
$ws.Cells.Item(9,14).Value = 571
$ws.Cells.Item(9,15).Value = 63.9224
$ws.Cells.Item(9,16).Value = 2.96598131277829
$ws.Cells.Item(9,17).Value = 5.7607095152339
$ws.Cells.Item(9,18).Value = 36
$ws.Cells.Item(9,19).Value = 0.36

Write-Output "edit complete"
